# Laborations_dagbok.xlsx update: "15/2-2018 Collision Test Stable"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new diary entry in row 3
$ws.Range("A3").Value = "15/2-2018"
$ws.Range("B3").Value = "Collision test"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 30

# Update the "Total tid" row so minutes overflow into hours
$ws.Range("C13").Formula = "=SUM(C2:C12)+QUOTIENT(SUM(D2:D12),60)"
$ws.Range("D13").Formula = "=MOD(SUM(D2:D12),60)"

# Move the active selection to D4, matching the author's final cursor position
$ws.Range("D4").Select()
